$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

# Column A holds a literal date-like text string (e.g. "12/08/2025"), matching
# the other rows in this sheet. Setting .Value directly would make Excel
# auto-convert the date-like text into a real date serial number, so instead
# we temporarily force a text number format, assign the text, then restore
# the default "Normal" style so no extra formatting lingers on the cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "12/08/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 13633.28
$ws.Cells.Item($row, 3).Value = 0.171641667611514
$ws.Cells.Item($row, 4).Value = 0.828358332388486
$ws.Cells.Item($row, 5).Value = -82.58
$ws.Cells.Item($row, 6).Value = -18.98
$ws.Cells.Item($row, 7).Value = -19198.52
$ws.Cells.Item($row, 8).Value = -62.96
$ws.Cells.Item($row, 9).Value = -462.81
$ws.Cells.Item($row, 10).Value = -16.51
